$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALXN")

# Row 4: Inventory
$ws.Range("C4").Value = 776000000.0
$ws.Range("D4").Value = 729000000.0
$ws.Range("E4").Value = 578000000.0
$ws.Range("F4").Value = 587000000.0
$ws.Range("G4").Value = 628000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 125300000.0
$ws.Range("C15").Value = 119000000.0
$ws.Range("D15").Value = 89000000.0
$ws.Range("E15").Value = 862000000.0
$ws.Range("F15").Value = 863000000.0
$ws.Range("G15").Value = 967000000.0

# Row 16: Accrued Expenses
$ws.Range("B16").Value = 910700000.0

# Row 24: Long Term Tax Liability (Deferred)
$ws.Range("C24").Value = -567000000.0
$ws.Range("D24").Value = -452000000.0
$ws.Range("E24").Value = -386000000.0
$ws.Range("F24").Value = -110000000.0
$ws.Range("G24").Value = -209000000.0

# Row 28: Additional Paid In Capital
$ws.Range("B28").Value = 9243300000.0

# Row 30: Retained Earnings
$ws.Range("B30").Value = 5879200000.0

# Row 31: Treasury Stock
$ws.Range("B31").Value = 2620500000.0

# Row 32: Common Equity (Total)
$ws.Range("B32").Value = 12431000000.0

# Row 35: Shares (Common)
$ws.Range("B35").Value = 220900000.0

# Row 37: Net Debt
$ws.Range("B37").Value = -937300000.0

# Row 38: Total Debt
$ws.Range("B38").Value = 2532000000.0
